# SectorGroup.xlsx: the "codeforiati:category-code" and
# "codeforiati:category-name" columns were swapped (category-name now
# comes before category-code), for the header row and every data row.
#
# We use Range.Copy(destination) rather than Value assignment so that
# numeric-looking text (e.g. category codes like "111") keeps its original
# text cell-type instead of being auto-coerced into a number, and so that
# no new cell styles / number formats get introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row()
$lastRow  = $used.Row() + $used.Rows.Count - 1
$firstCol = $used.Column()
$lastCol  = $used.Column() + $used.Columns.Count - 1

# Locate the two columns to swap by their header text (row 1).
$codeCol = -1
$nameCol = -1
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Value()
    if ($header -eq "codeforiati:category-code") { $codeCol = $c }
    if ($header -eq "codeforiati:category-name") { $nameCol = $c }
}

if ($codeCol -gt 0 -and $nameCol -gt 0) {
    # Pick a scratch column well clear of the used range so it can't
    # collide with real data and won't survive in the saved dimension.
    $tempCol = $lastCol + 2

    $codeRange = $ws.Range($ws.Cells.Item($firstRow, $codeCol), $ws.Cells.Item($lastRow, $codeCol))
    $nameRange = $ws.Range($ws.Cells.Item($firstRow, $nameCol), $ws.Cells.Item($lastRow, $nameCol))
    $tempRange = $ws.Range($ws.Cells.Item($firstRow, $tempCol), $ws.Cells.Item($lastRow, $tempCol))

    $codeRange.Copy($tempRange)
    $nameRange.Copy($codeRange)
    $tempRange.Copy($nameRange)

    $tempRange.Clear()
}
